$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.072.15'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '1.646.15'
$ws.Range("E3").Value = '  +0.74%  '
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.63'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.65'
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("D12").Value = '1.876.12'
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.29'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = '1.672.33'
$ws.Range("E14").Value = '  +3.17%  '
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.45'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '26.310.86'
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.61'
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.94'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.81'
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("E25").Value = '  +4.00%  '
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '143.94'
$ws.Range("E27").Value = '  +0.67%  '
$ws.Range("E28").Value = '  +0.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.56'
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0497'
$ws.Range("E32").Value = '  +1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.30'
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("D37").Value = '1.134.11'
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.542'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.46'
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D44").Value = '1.784.74'
$ws.Range("E44").Value = '  +0.82%  '
$ws.Range("D45").Value = '0.0₆0116'
$ws.Range("E45").Value = '  +4.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.72'
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.76'
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("E51").Value = '  +0.09%  '
